$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 146.375
$ws.Range("I12").Value = 125
$ws.Range("J12").Value = 182
$ws.Range("K12").Value = 125
$ws.Range("L12").Value = 182
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = -522
$ws.Range("H28").Value = 634.6667
$ws.Range("J28").Value = 1236.4
$ws.Range("L28").Value = 1236.4
$ws.Range("N28").Value = -2206.4
$ws.Range("H86").Value = 4682.4585
$ws.Range("I86").Value = 3641.8572
$ws.Range("K86").Value = 3641.8572
$ws.Range("M86").Value = -2518.8572
$ws.Range("H89").Value = 4682.4585
$ws.Range("I89").Value = 3641.8572
$ws.Range("K89").Value = 18209.286
$ws.Range("M89").Value = -12593.286
$ws.Range("H116").Value = 210867.3
$ws.Range("I116").Value = 49384.848
$ws.Range("K116").Value = 49384.848
$ws.Range("M116").Value = -45942.848
$ws.Range("H135").Value = 658.6316
$ws.Range("I135").Value = 546.625
$ws.Range("J135").Value = 1256
$ws.Range("K135").Value = 4919.625
$ws.Range("L135").Value = 11304
$ws.Range("M135").Value = -2384.625
$ws.Range("N135").Value = -16374
$ws.Range("H137").Value = 13079.5
$ws.Range("I137").Value = 2161.25
$ws.Range("J137").Value = 23997.75
$ws.Range("K137").Value = 6483.75
$ws.Range("L137").Value = 71993.25
$ws.Range("M137").Value = -3933.75
$ws.Range("N137").Value = -77093.25
$ws.Range("H138").Value = 2920.3
$ws.Range("J138").Value = 4860.409
$ws.Range("L138").Value = 14581.227
$ws.Range("N138").Value = -24861.227

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16162.571
$ws.Range("I2").Value = 6059.8
$ws.Range("J2").Value = 21775.223
$ws.Range("K2").Value = 6059.8
$ws.Range("L2").Value = 21775.223
$ws.Range("M2").Value = -5946.8
$ws.Range("N2").Value = -22001.223
$ws.Range("H61").Value = 879619.3
$ws.Range("I61").Value = 982813.25
$ws.Range("J61").Value = 2471
$ws.Range("K61").Value = 982813.25
$ws.Range("L61").Value = 2471
$ws.Range("M61").Value = -982601.25
$ws.Range("N61").Value = -2895
$ws.Range("H63").Value = 15000.071
$ws.Range("I63").Value = 4498.5
$ws.Range("K63").Value = 4498.5
$ws.Range("M63").Value = -3812.5
$ws.Range("H66").Value = 15000.071
$ws.Range("I66").Value = 4498.5
$ws.Range("K66").Value = 22492.5
$ws.Range("M66").Value = -19060.5
$ws.Range("H74").Value = 3782.8
$ws.Range("I74").Value = 2134.8572
$ws.Range("K74").Value = 2134.8572
$ws.Range("M74").Value = -1260.8572
$ws.Range("H77").Value = 3782.8
$ws.Range("I77").Value = 2134.8572
$ws.Range("K77").Value = 10674.286
$ws.Range("M77").Value = -6306.286
$ws.Range("H101").Value = 49599.8
$ws.Range("J101").Value = 49599.8
$ws.Range("L101").Value = 49599.8
$ws.Range("N101").Value = -56089.8
$ws.Range("H116").Value = 16162.571
$ws.Range("I116").Value = 6059.8
$ws.Range("J116").Value = 21775.223
$ws.Range("K116").Value = 6059.8
$ws.Range("L116").Value = 21775.223
$ws.Range("M116").Value = -3765.8
$ws.Range("N116").Value = -26363.223
$ws.Range("H122").Value = 2215.5833
$ws.Range("I122").Value = 1632.625
$ws.Range("J122").Value = 3381.5
$ws.Range("K122").Value = 4897.875
$ws.Range("L122").Value = 10144.5
$ws.Range("M122").Value = -2447.875
$ws.Range("N122").Value = -15044.5
$ws.Range("H132").Value = 1054829.2
$ws.Range("I132").Value = 1251734.8
$ws.Range("K132").Value = 3755204.4
$ws.Range("M132").Value = -3752674.4
$ws.Range("H136").Value = 879619.3
$ws.Range("I136").Value = 982813.25
$ws.Range("J136").Value = 2471
$ws.Range("K136").Value = 2948439.75
$ws.Range("L136").Value = 7413
$ws.Range("M136").Value = -2945889.75
$ws.Range("N136").Value = -12513

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16162.571
$ws.Range("I3").Value = 6059.8
$ws.Range("J3").Value = 21775.223
$ws.Range("K3").Value = 6059.8
$ws.Range("L3").Value = 21775.223
$ws.Range("M3").Value = -5945.8
$ws.Range("N3").Value = -22003.223
$ws.Range("H86").Value = 1127.1111
$ws.Range("I86").Value = 1381.3334
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1381.3334
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -258.3334
$ws.Range("N86").Value = -3246
$ws.Range("H89").Value = 1127.1111
$ws.Range("I89").Value = 1381.3334
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6906.666999999999
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -1290.666999999999
$ws.Range("N89").Value = -16232
$ws.Range("H99").Value = 18445
$ws.Range("I99").Value = 22603
$ws.Range("K99").Value = 22603
$ws.Range("M99").Value = -21105
$ws.Range("H134").Value = 980441.5600000001
$ws.Range("I134").Value = 993641.7
$ws.Range("J134").Value = 917081
$ws.Range("K134").Value = 2980925.1
$ws.Range("L134").Value = 2751243
$ws.Range("M134").Value = -2978390.1
$ws.Range("N134").Value = -2756313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 689710.5
$ws.Range("I58").Value = 1125407.1
$ws.Range("K58").Value = 1125407.1
$ws.Range("M58").Value = -1125204.1
$ws.Range("H93").Value = 36599.668
$ws.Range("I93").Value = 39899.5
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 39899.5
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -38027.5
$ws.Range("N93").Value = -33744
$ws.Range("H132").Value = 210617
$ws.Range("I132").Value = 1806
$ws.Range("J132").Value = 837050
$ws.Range("K132").Value = 5418
$ws.Range("L132").Value = 2511150
$ws.Range("M132").Value = -2888
$ws.Range("N132").Value = -2516210
$ws.Range("H134").Value = 1695.9286
$ws.Range("I134").Value = 1495.3334
$ws.Range("K134").Value = 4486.0002
$ws.Range("M134").Value = -1951.0002
$ws.Range("H136").Value = 689710.5
$ws.Range("I136").Value = 1125407.1
$ws.Range("K136").Value = 3376221.3
$ws.Range("M136").Value = -3373671.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22848650
$ws.Range("I4").Value = 33036124
$ws.Range("K4").Value = 99108372
$ws.Range("M4").Value = -99108260
$ws.Range("H14").Value = 157.9
$ws.Range("I14").Value = 157.9
$ws.Range("K14").Value = 473.7
$ws.Range("M14").Value = -300.7
$ws.Range("H40").Value = 88.125
$ws.Range("I40").Value = 29.428572
$ws.Range("J40").Value = 499
$ws.Range("K40").Value = 117.714288
$ws.Range("L40").Value = 1996
$ws.Range("M40").Value = -48.714288
$ws.Range("N40").Value = -2134
$ws.Range("H69").Value = 3073.75
$ws.Range("I69").Value = 1150
$ws.Range("J69").Value = 4997.5
$ws.Range("K69").Value = 3450
$ws.Range("L69").Value = 14992.5
$ws.Range("M69").Value = -2639
$ws.Range("N69").Value = -16614.5
$ws.Range("H72").Value = 3073.75
$ws.Range("I72").Value = 1150
$ws.Range("J72").Value = 4997.5
$ws.Range("K72").Value = 10350
$ws.Range("L72").Value = 44977.5
$ws.Range("M72").Value = -6294
$ws.Range("N72").Value = -53089.5
$ws.Range("H107").Value = 659.4375
$ws.Range("J107").Value = 603
$ws.Range("L107").Value = 1809
$ws.Range("N107").Value = -5649
$ws.Range("H140").Value = 2673.5417
$ws.Range("I140").Value = 1464.1666
$ws.Range("K140").Value = 4392.4998
$ws.Range("M140").Value = 787.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12662.5
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("N5").Value = -324
$ws.Range("H45").Value = 46300
$ws.Range("J45").Value = 46300
$ws.Range("L45").Value = 46300
$ws.Range("N45").Value = -47418
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 10666.714
$ws.Range("I132").Value = 9648.556
$ws.Range("J132").Value = 12499.4
$ws.Range("K132").Value = 28945.668
$ws.Range("L132").Value = 37498.2
$ws.Range("M132").Value = -26415.668
$ws.Range("N132").Value = -42558.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1631.5
$ws.Range("J22").Value = 2499.6667
$ws.Range("L22").Value = 2499.6667
$ws.Range("N22").Value = -3089.6667
$ws.Range("H27").Value = 1631.5
$ws.Range("J27").Value = 2499.6667
$ws.Range("L27").Value = 2499.6667
$ws.Range("N27").Value = -2713.6667
$ws.Range("H122").Value = 4574.643
$ws.Range("I122").Value = 4171.0415
$ws.Range("J122").Value = 6996.25
$ws.Range("K122").Value = 12513.1245
$ws.Range("L122").Value = 20988.75
$ws.Range("M122").Value = -10063.1245
$ws.Range("N122").Value = -25888.75
$ws.Range("H132").Value = 25401.637
$ws.Range("I132").Value = 29657.555
$ws.Range("K132").Value = 88972.66500000001
$ws.Range("M132").Value = -86442.66500000001
$ws.Range("H136").Value = 2998.6296
$ws.Range("I136").Value = 2123.55
$ws.Range("K136").Value = 6370.650000000001
$ws.Range("M136").Value = -3820.650000000001
$ws.Range("H139").Value = 92357.5
$ws.Range("J139").Value = 92357.5
$ws.Range("L139").Value = 92357.5
$ws.Range("N139").Value = -102637.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1910.0385
$ws.Range("I122").Value = 1582.4736
$ws.Range("J122").Value = 2799.1428
$ws.Range("K122").Value = 4747.4208
$ws.Range("L122").Value = 8397.428400000001
$ws.Range("M122").Value = -2297.4208
$ws.Range("N122").Value = -13297.4284
$ws.Range("H132").Value = 2827269.8
$ws.Range("I132").Value = 3207383.5
$ws.Range("K132").Value = 9622150.5
$ws.Range("M132").Value = -9619620.5
$ws.Range("H136").Value = 6675.2144
$ws.Range("I136").Value = 6791.057
$ws.Range("J136").Value = 6096
$ws.Range("K136").Value = 20373.171
$ws.Range("L136").Value = 18288
$ws.Range("M136").Value = -17823.171
